$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Region", "Species", "LifeStage", "Season"),
    @("North Delta", "Giant Reed EAV", "", "Spring"),
    @("North Delta", "Giant Reed EAV", "", "Summer"),
    @("North Delta", "Giant Reed EAV", "", "Fall"),
    @("North Delta", "Giant Reed EAV", "", "Winter"),
    @("South", "Giant Reed EAV", "", "Spring"),
    @("South", "Giant Reed EAV", "", "Summer"),
    @("South", "Giant Reed EAV", "", "Fall"),
    @("South", "Giant Reed EAV", "", "Winter"),
    @("Central", "Giant Reed EAV", "", "Spring"),
    @("Central", "Giant Reed EAV", "", "Summer"),
    @("Central", "Giant Reed EAV", "", "Fall"),
    @("Central", "Giant Reed EAV", "", "Winter"),
    @("Suisun Bay", "Giant Reed EAV", "", "Spring"),
    @("Suisun Bay", "Giant Reed EAV", "", "Summer"),
    @("Suisun Bay", "Giant Reed EAV", "", "Fall"),
    @("Suisun Bay", "Giant Reed EAV", "", "Winter"),
    @("Suisun Marsh", "Giant Reed EAV", "", "Spring"),
    @("Suisun Marsh", "Giant Reed EAV", "", "Summer"),
    @("Suisun Marsh", "Giant Reed EAV", "", "Fall"),
    @("Suisun Marsh", "Giant Reed EAV", "", "Winter"),
    @("Confluence", "Giant Reed EAV", "", "Spring"),
    @("Confluence", "Giant Reed EAV", "", "Summer"),
    @("Confluence", "Giant Reed EAV", "", "Fall"),
    @("Confluence", "Giant Reed EAV", "", "Winter"),
    @("North Delta", "Brazilian Waterweed", "", "Spring"),
    @("North Delta", "Brazilian Waterweed", "", "Summer"),
    @("North Delta", "Brazilian Waterweed", "", "Fall"),
    @("North Delta", "Brazilian Waterweed", "", "Winter"),
    @("South", "Brazilian Waterweed", "", "Spring"),
    @("South", "Brazilian Waterweed", "", "Summer"),
    @("South", "Brazilian Waterweed", "", "Fall"),
    @("South", "Brazilian Waterweed", "", "Winter"),
    @("Central", "Brazilian Waterweed", "", "Spring"),
    @("Central", "Brazilian Waterweed", "", "Summer"),
    @("Central", "Brazilian Waterweed", "", "Fall"),
    @("Central", "Brazilian Waterweed", "", "Winter"),
    @("Suisun Bay", "Brazilian Waterweed", "", "Spring"),
    @("Suisun Bay", "Brazilian Waterweed", "", "Summer"),
    @("Suisun Bay", "Brazilian Waterweed", "", "Fall"),
    @("Suisun Bay", "Brazilian Waterweed", "", "Winter"),
    @("Suisun Marsh", "Brazilian Waterweed", "", "Spring"),
    @("Suisun Marsh", "Brazilian Waterweed", "", "Summer"),
    @("Suisun Marsh", "Brazilian Waterweed", "", "Fall"),
    @("Suisun Marsh", "Brazilian Waterweed", "", "Winter"),
    @("Confluence", "Brazilian Waterweed", "", "Spring"),
    @("Confluence", "Brazilian Waterweed", "", "Summer"),
    @("Confluence", "Brazilian Waterweed", "", "Fall"),
    @("Confluence", "Brazilian Waterweed", "", "Winter"),
    @("North Delta", "Water Hyacinth FAV", "", "Spring"),
    @("North Delta", "Water Hyacinth FAV", "", "Summer"),
    @("North Delta", "Water Hyacinth FAV", "", "Fall"),
    @("North Delta", "Water Hyacinth FAV", "", "Winter"),
    @("South", "Water Hyacinth FAV", "", "Spring"),
    @("South", "Water Hyacinth FAV", "", "Summer"),
    @("South", "Water Hyacinth FAV", "", "Fall"),
    @("South", "Water Hyacinth FAV", "", "Winter"),
    @("Central", "Water Hyacinth FAV", "", "Spring"),
    @("Central", "Water Hyacinth FAV", "", "Summer"),
    @("Central", "Water Hyacinth FAV", "", "Fall"),
    @("Central", "Water Hyacinth FAV", "", "Winter"),
    @("Suisun Bay", "Water Hyacinth FAV", "", "Spring"),
    @("Suisun Bay", "Water Hyacinth FAV", "", "Summer"),
    @("Suisun Bay", "Water Hyacinth FAV", "", "Fall"),
    @("Suisun Bay", "Water Hyacinth FAV", "", "Winter"),
    @("Suisun Marsh", "Water Hyacinth FAV", "", "Spring"),
    @("Suisun Marsh", "Water Hyacinth FAV", "", "Summer"),
    @("Suisun Marsh", "Water Hyacinth FAV", "", "Fall"),
    @("Suisun Marsh", "Water Hyacinth FAV", "", "Winter"),
    @("Confluence", "Water Hyacinth FAV", "", "Spring"),
    @("Confluence", "Water Hyacinth FAV", "", "Summer"),
    @("Confluence", "Water Hyacinth FAV", "", "Fall"),
    @("Confluence", "Water Hyacinth FAV", "", "Winter"),
    @("North Delta", "Microcystis", "", "Summer"),
    @("North Delta", "Microcystis", "", "Fall"),
    @("South", "Microcystis", "", "Summer"),
    @("South", "Microcystis", "", "Fall"),
    @("Central", "Microcystis", "", "Summer"),
    @("Central", "Microcystis", "", "Fall"),
    @("Suisun Bay", "Microcystis", "", "Summer"),
    @("Suisun Bay", "Microcystis", "", "Fall"),
    @("Suisun Marsh", "Microcystis", "", "Summer"),
    @("Suisun Marsh", "Microcystis", "", "Fall"),
    @("Confluence", "Microcystis", "", "Summer"),
    @("Confluence", "Microcystis", "", "Fall"),
    @("Suisun Bay", "Overbite Clam", "", "Summer"),
    @("Suisun Bay", "Overbite Clam", "", "Fall"),
    @("Suisun Bay", "Overbite Clam", "", "Spring"),
    @("Suisun Bay", "Overbite Clam", "", "Winter"),
    @("Suisun Marsh", "Overbite Clam", "", "Summer"),
    @("Suisun Marsh", "Overbite Clam", "", "Fall"),
    @("Suisun Marsh", "Overbite Clam", "", "Spring"),
    @("Suisun Marsh", "Overbite Clam", "", "Winter"),
    @("Confluence", "Overbite Clam", "", "Summer"),
    @("Confluence", "Overbite Clam", "", "Fall"),
    @("Confluence", "Overbite Clam", "", "Spring"),
    @("Confluence", "Overbite Clam", "", "Winter"),
    @("North Delta", "Asian Clam", "", "Spring"),
    @("North Delta", "Asian Clam", "", "Summer"),
    @("North Delta", "Asian Clam", "", "Fall"),
    @("North Delta", "Asian Clam", "", "Winter"),
    @("South", "Asian Clam", "", "Spring"),
    @("South", "Asian Clam", "", "Summer"),
    @("South", "Asian Clam", "", "Fall"),
    @("South", "Asian Clam", "", "Winter"),
    @("Central", "Asian Clam", "", "Spring"),
    @("Central", "Asian Clam", "", "Summer"),
    @("Central", "Asian Clam", "", "Fall"),
    @("Central", "Asian Clam", "", "Winter"),
    @("Suisun Bay", "Asian Clam", "", "Spring"),
    @("Suisun Bay", "Asian Clam", "", "Summer"),
    @("Suisun Bay", "Asian Clam", "", "Fall"),
    @("Suisun Bay", "Asian Clam", "", "Winter"),
    @("Suisun Marsh", "Asian Clam", "", "Spring"),
    @("Suisun Marsh", "Asian Clam", "", "Summer"),
    @("Suisun Marsh", "Asian Clam", "", "Fall"),
    @("Suisun Marsh", "Asian Clam", "", "Winter"),
    @("Confluence", "Asian Clam", "", "Spring"),
    @("Confluence", "Asian Clam", "", "Summer"),
    @("Confluence", "Asian Clam", "", "Fall"),
    @("Confluence", "Asian Clam", "", "Winter")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r+1, $c+1).Value = $row[$c]
    }
}

# Microcystis rows for Suisun Bay/Suisun Marsh/Confluence (rows 80-85) use the
# same italic species-column styling as the other Microcystis rows (74-79).
$ws.Range("B80:B85").Font.Italic = $true

# Header row / sheet selection tidy-up to match the saved view state.
[void]$ws.Range("C2").Select()
